# This script applies the "Updated symbol list on Fri Feb 17 09:35:23 UTC 2023
# with GitHub Actions" edit to the cryptocurrency tracking sheet.
#
# All data cells in columns D (Price) and E (Volume/1h) are stored as plain
# text strings (e.g. "309.92", "-3.63%") rather than numeric values, so we
# must make sure Excel keeps them as text instead of silently re-interpreting
# them as numbers/percentages. We do this by prefixing the value with a
# leading apostrophe (the normal Excel "treat as text" quote-prefix) and then
# resetting the cell style back to "Normal" so no stray number-format / style
# is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

# Row 2 (D,E)
Set-TextValue "D2" '310.13'
Set-TextValue "E2" '-3.73%'

# Row 3 (D,E)
Set-TextValue "D3" '49.50'
Set-TextValue "E3" '2.03%'

# Row 4 (D,E)
Set-TextValue "D4" '5.132'
Set-TextValue "E4" '-2.51%'

# Row 5 (D,E)
Set-TextValue "D5" '0.07769'
Set-TextValue "E5" '-4.21%'

# Row 6 (D,E)
Set-TextValue "D6" '4.537'
Set-TextValue "E6" '-0.98%'

# Row 7 (D,E)
Set-TextValue "D7" '1.379'
Set-TextValue "E7" '14.91%'

# Row 8 (D,E)
Set-TextValue "D8" '1.575'
Set-TextValue "E8" '-4.14%'

# Row 9 (D,E)
Set-TextValue "D9" '0.1218'
Set-TextValue "E9" '-6.23%'

# Row 10 (E)
Set-TextValue "E10" '1.70%'

# Row 11 (D,E)
Set-TextValue "D11" '0.04701'
Set-TextValue "E11" '1.54%'

# Row 12 (E)
Set-TextValue "E12" '-2.13%'

# Row 13 (E)
Set-TextValue "E13" '-0.48%'

# Row 14 (D,E)
Set-TextValue "D14" '0.001275'
Set-TextValue "E14" '-3.83%'

# Row 15 (D,E)
Set-TextValue "D15" '0.04171'
Set-TextValue "E15" '-0.68%'

# Row 16 (D,E)
Set-TextValue "D16" '0.005791'
Set-TextValue "E16" '-0.97%'

# Row 17 (E)
Set-TextValue "E17" '2,021.10%'

# Row 18 (E)
Set-TextValue "E18" '-0.18%'

# Row 19 (D,E)
Set-TextValue "D19" '2.434'
Set-TextValue "E19" '0.32%'

# Row 20 (D,E)
Set-TextValue "D20" '0.3396'
Set-TextValue "E20" '-0.16%'

# Row 21 (D,E)
Set-TextValue "D21" '7.983'
Set-TextValue "E21" '-1.44%'

# Row 22 (D,E)
Set-TextValue "D22" '0.1337'
Set-TextValue "E22" '-5.21%'

# Row 23 (D,E)
Set-TextValue "D23" '0.3036'
Set-TextValue "E23" '-2.85%'

# Row 24 (D,E)
Set-TextValue "D24" '0.001271'
Set-TextValue "E24" '-2.68%'

# Row 25 (D,E)
Set-TextValue "D25" '0.003936'
Set-TextValue "E25" '-7.42%'

# Row 26 (D,E)
Set-TextValue "D26" '0.0001349'
Set-TextValue "E26" '-0.06%'

# Row 38 (D,E)
Set-TextValue "D38" '0.02596'
Set-TextValue "E38" '-3.67%'

# Row 39 (D,E)
Set-TextValue "D39" '0.06256'
Set-TextValue "E39" '10.73%'

# Row 41 (D,E)
Set-TextValue "D41" '0.007934'
Set-TextValue "E41" '3.36%'

# Row 42 (D,E)
Set-TextValue "D42" '0.1422'
Set-TextValue "E42" '-1.24%'

# Row 43 (D,E)
Set-TextValue "D43" '0.008373'
Set-TextValue "E43" '8.79%'

# Row 44 (D,E)
Set-TextValue "D44" '0.008334'
Set-TextValue "E44" '2.87%'

# Row 45 (D,E)
Set-TextValue "D45" '0.3129'
Set-TextValue "E45" '-2.03%'

# Row 46 (D,E)
Set-TextValue "D46" '0.00007674'
Set-TextValue "E46" '9.52%'

# Row 47 (D,E)
Set-TextValue "D47" '0.00000000750'
Set-TextValue "E47" '-0.06%'

# Row 48 (B,C,D,E)
Set-TextValue "B48" 'CoinbaseStockToken'
Set-TextValue "C48" 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue "D48" '0.002618'
Set-TextValue "E48" '-34.56%'

# Row 49 (B,C,D,E)
Set-TextValue "B49" 'BOLO'
Set-TextValue "C49" 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue "D49" '0.05321'
Set-TextValue "E49" '-1.56%'

# Row 50 (D,E)
Set-TextValue "D50" '0.00002099'
Set-TextValue "E50" '-0.06%'

# Row 51 (D,E)
Set-TextValue "D51" '0.0001999'
Set-TextValue "E51" '-0.06%'
